$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. Insert a fresh row
# at position 33 (pushing the existing rows 33-117 down to 34-118) and
# populate it with the new observation.
$ws.Rows.Item(33).Insert()

$ws.Range("A33").Value = 4
$ws.Range("B33").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C33").Value = "Los Lagos"
$ws.Range("D33").Value = 44498
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = 100112009
$ws.Range("G33").Value = "Acelga"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = 3000
$ws.Range("N33").Value = "$/docena de atados (4 kilos)"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 750
$ws.Range("Q33").Value = 4
$ws.Range("R33").Value = "Hortaliza"
